$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D7 previously held the label "B7:C8" - clear it out.
$ws.Range("D7").Value = ""

# E7: "Privet" with the middle two letters ("ri") in italics.
$ws.Range("E7").Value = "Privet"
$ws.Range("E7").Characters(2, 2).Font.Italic = $true

# G7: "fontname" with the second half ("name") underlined (double),
# colored blue and sized 20.
$ws.Range("G7").Value = "fontname"
$nameChars = $ws.Range("G7").Characters(5, 4)
$nameChars.Font.Underline = $true
$nameChars.Font.Color = 16711680
$nameChars.Font.Size = 20
